# Updated BGR model - 2025-08-08 15:56
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEDA_Sets-Proc")

# Row 19 defines the "Util Batt Stg" processor set. Tighten its PSET_PN
# match expression (column B) from a generic "*bat*" wildcard to the more
# specific "EN*STG?hb*" pattern (still excluding "*EV*").
$ws.Range("B19").Value = "EN*STG?hb*,-*EV*"

# Also record the And/Or operators for this set's positive/negative
# wildcard lists (columns H / I), matching the pattern already used on
# other rows (e.g. row 17) that combine multiple match terms.
$ws.Range("H19").Value = "And"
$ws.Range("I19").Value = "Or"
